# Buffer calculation function and point-data processing script
#
# Adds the FEMA National Risk Index dataset as a new row to the
# indicator_search tracking sheet, and leaves the sheet scrolled/selected
# near the bottom of the table (mirrors where the author was working).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19
$ws.Cells.Item($row, 1).Value = "National Risk Dataset"
$ws.Cells.Item($row, 2).Value = "CONUS"
$ws.Cells.Item($row, 3).Value = "Y"
$ws.Cells.Item($row, 4).Value = ".csv, .shp, .gdb"
$ws.Cells.Item($row, 5).Value = "Census tract level"
$ws.Cells.Item($row, 6).Value = "~"
$ws.Cells.Item($row, 7).Value = "FEMA"
$ws.Cells.Item($row, 8).Value = "https://hazards.fema.gov/nri/data-resources"

# Match the author's final on-screen selection after entering the row.
$ws.Range("D20").Select()
